$p = $ppt.ActivePresentation
$dt = $p.DocumentTheme
$tfs = $dt.ThemeFontScheme
$mf = $tfs.MajorFont
Get-Member -InputObject $mf
